# Adds a new "دی 99" monthly timesheet block (rows 139-151) to Sheet1,
# mirroring the structure of the previous month's block (rows 126-138).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

function Copy-Format($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats) | Out-Null
}

# --- Row 139: blank thick-bottom separator row (mirrors row 126) ---
Copy-Format "A126:G126" "A139:G139"
$ws.Rows(139).RowHeight = 15

# --- Row 140: blank shaded row (mirrors row 127) ---
Copy-Format "A127:E127" "A140:E140"

# --- Row 141: month header row (mirrors row 128) ---
Copy-Format "A128:E128" "A141:E141"
$ws.Range("A141").Value = "دی 99"
$ws.Range("B141").Value = "Activity"
$ws.Range("C141").Value = "Hours"
$ws.Range("E141").Value = "Tasks Done"

# --- Row 142: (mirrors row 129) ---
Copy-Format "B129" "B142"
Copy-Format "C129" "C142"
Copy-Format "E129" "E142"
$ws.Range("B142").Value = "* Documentation"
$ws.Range("E142").Value = [char]0x2022 + " "

# --- Row 143: (mirrors row 130) ---
Copy-Format "B130" "B143"
Copy-Format "C130" "C143"
Copy-Format "E130" "E143"
$ws.Range("B143").Value = "* GUI"
$ws.Range("E143").Value = [char]0x2022 + " "

# --- Row 144: (mirrors row 131) ---
Copy-Format "B131" "B144"
Copy-Format "C131" "C144"
Copy-Format "E131" "E144"
$ws.Range("B144").Value = "* Registration"
$ws.Range("E144").Value = [char]0x2022 + " "

# --- Row 145: (mirrors row 132) ---
Copy-Format "B132" "B145"
Copy-Format "C132" "C145"
Copy-Format "E132" "E145"
$ws.Range("B145").Value = "* Tracker"
$ws.Range("E145").Value = [char]0x2022 + " "

# --- Row 146: (mirrors row 133) ---
Copy-Format "B133" "B146"
Copy-Format "C133" "C146"
Copy-Format "E133" "E146"
$ws.Range("B146").Value = "* 2D/3D Views"
$ws.Range("E146").Value = [char]0x2022 + " "

# --- Row 147: (mirrors row 134, no E cell) ---
Copy-Format "B134" "B147"
Copy-Format "C134" "C147"
$ws.Range("B147").Value = "* Patients / Database"

# --- Row 148: Paid hours row with payment formula (mirrors row 135) ---
Copy-Format "B135" "B148"
Copy-Format "C135" "C148"
Copy-Format "E135" "E148"
$ws.Rows(148).RowHeight = 15.6
$ws.Range("B148").Value = "* Meetings & Presentations"
$ws.Range("C148").Value = 1
$ws.Range("E148").Formula = "=C149*40000"

# --- Row 149: Total hours row (mirrors row 136) ---
Copy-Format "B136" "B149"
Copy-Format "C136" "C149"
$ws.Range("B149").Value = [char]0x2022 + " Total Hours"
$ws.Range("C149").Formula = "=SUM(C142:C148)"

# --- Row 150: Paid total (mirrors row 137) ---
Copy-Format "C137" "C150"
Copy-Format "D137" "D150"
$ws.Range("C150").Value = "@Parsiss"
$ws.Range("D150").Value = 0

# --- Row 151: Not-paid total (mirrors row 138) ---
Copy-Format "C138" "C151"
Copy-Format "D138" "D151"
$ws.Range("C151").Value = "@Home"
$ws.Range("D151").Formula = "=C149-D150"

# --- Update the view: scroll/selection to match the new content ---
$ws.Range("E146").Select()

Write-Host "Added December (دی 99) timesheet block: rows 139-151"
